$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 336, shifting rows 336:392 down to 337:393
$ws.Rows.Item(336).Insert()

# Populate the newly inserted row 336 with the new data record
$ws.Range("A336").Value = 9
$ws.Range("B336").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C336").Value = "Metropolitana"
$ws.Range("D336").Value = 44637
$ws.Range("E336").Value = 13
$ws.Range("F336").Value = 100112012
$ws.Range("G336").Value = "Espinaca"
$ws.Range("H336").Value = "Sin especificar"
$ws.Range("I336").Value = "Primera"
$ws.Range("J336").Value = 160
$ws.Range("K336").Value = 10000
$ws.Range("L336").Value = 12000
$ws.Range("M336").Value = 11000
$ws.Range("N336").Value = "`$/cuna 10 kilos"
$ws.Range("O336").Value = "Provincia de Chacabuco"
$ws.Range("P336").Value = 1100
$ws.Range("Q336").Value = 10
$ws.Range("R336").Value = "Hortaliza"

# Match the date cell style used by the other date cells in column D (style index 2)
$ws.Range("D336").NumberFormat = $ws.Range("D337").NumberFormat
